# Update "想去人数" (number of people interested) counts for three rows.
# These values live in column F, rows 3-5, and must be updated identically
# on both the "展览" sheet and the "全部类型" sheet (they mirror each other).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 219
    $ws.Range("F4").Value = 3714
    $ws.Range("F5").Value = 383
}
